$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the rate for Navel from 140 to 400
$ws.Range("B2").Value = 400
